$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4734.6924
$ws.Range("I64").Value = 4999
$ws.Range("J64").Value = 4724.12
$ws.Range("K64").Value = 4999
$ws.Range("L64").Value = 4724.12
$ws.Range("M64").Value = -4751
$ws.Range("N64").Value = -5220.12
$ws.Range("H67").Value = 4734.6924
$ws.Range("I67").Value = 4999
$ws.Range("J67").Value = 4724.12
$ws.Range("K67").Value = 4999
$ws.Range("L67").Value = 4724.12
$ws.Range("M67").Value = -4141
$ws.Range("N67").Value = -6440.12
$ws.Range("H76").Value = 4254.4614
$ws.Range("I76").Value = 4215.4287
$ws.Range("J76").Value = 4300
$ws.Range("K76").Value = 4215.4287
$ws.Range("L76").Value = 4300
$ws.Range("M76").Value = -3900.4287
$ws.Range("N76").Value = -4930
$ws.Range("H79").Value = 4254.4614
$ws.Range("I79").Value = 4215.4287
$ws.Range("J79").Value = 4300
$ws.Range("K79").Value = 4215.4287
$ws.Range("L79").Value = 4300
$ws.Range("M79").Value = -3123.4287
$ws.Range("N79").Value = -6484
$ws.Range("H98").Value = 6137.8066
$ws.Range("I98").Value = 4717.885
$ws.Range("J98").Value = 13521.4
$ws.Range("K98").Value = 4717.885
$ws.Range("L98").Value = 13521.4
$ws.Range("M98").Value = -3219.885
$ws.Range("N98").Value = -16517.4
$ws.Range("H107").Value = 125749.375
$ws.Range("I107").Value = 143585
$ws.Range("J107").Value = 900
$ws.Range("K107").Value = 143585
$ws.Range("L107").Value = 900
$ws.Range("M107").Value = -141665
$ws.Range("N107").Value = -4740
$ws.Range("H122").Value = 6137.8066
$ws.Range("I122").Value = 4717.885
$ws.Range("J122").Value = 13521.4
$ws.Range("K122").Value = 14153.655
$ws.Range("L122").Value = 40564.2
$ws.Range("M122").Value = -11703.655
$ws.Range("N122").Value = -45464.2
$ws.Range("H129").Value = 1000.0263
$ws.Range("J129").Value = 1234
$ws.Range("L129").Value = 3702
$ws.Range("N129").Value = -13702

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H38").Value = 10000
$ws.Range("J38").Value = 10000
$ws.Range("L38").Value = 10000
$ws.Range("N38").Value = -10934
$ws.Range("H122").Value = 2059.2917
$ws.Range("I122").Value = 1801.15
$ws.Range("J122").Value = 3350
$ws.Range("K122").Value = 5403.450000000001
$ws.Range("L122").Value = 10050
$ws.Range("M122").Value = -2953.450000000001
$ws.Range("N122").Value = -14950

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3596.5454
$ws.Range("I105").Value = 3596.5454
$ws.Range("K105").Value = 3596.5454
$ws.Range("M105").Value = -1849.5454

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3116.6316
$ws.Range("I31").Value = 1688.9688
$ws.Range("K31").Value = 1688.9688
$ws.Range("M31").Value = -1393.9688
$ws.Range("H34").Value = 3116.6316
$ws.Range("I34").Value = 1688.9688
$ws.Range("K34").Value = 1688.9688
$ws.Range("M34").Value = -1486.9688
$ws.Range("H123").Value = 60780
$ws.Range("J123").Value = 60780
$ws.Range("L123").Value = 60780
$ws.Range("N123").Value = -70580

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 371.9
$ws.Range("I23").Value = 243.33333
$ws.Range("J23").Value = 427
$ws.Range("K23").Value = 729.99999
$ws.Range("L23").Value = 1281
$ws.Range("M23").Value = -494.99999
$ws.Range("N23").Value = -1751
$ws.Range("H96").Value = 27100
$ws.Range("J96").Value = 4200
$ws.Range("L96").Value = 12600
$ws.Range("N96").Value = -16718
$ws.Range("H112").Value = 3517.75
$ws.Range("I112").Value = 1205.375
$ws.Range("J112").Value = 4288.5415
$ws.Range("K112").Value = 3616.125
$ws.Range("L112").Value = 12865.6245
$ws.Range("M112").Value = -2508.125
$ws.Range("N112").Value = -15081.6245
$ws.Range("H131").Value = 856.28
$ws.Range("J131").Value = 879.55316
$ws.Range("L131").Value = 2638.65948
$ws.Range("N131").Value = -12718.65948
$ws.Range("H132").Value = 2121.9
$ws.Range("I132").Value = 1484.2222
$ws.Range("J132").Value = 2643.6365
$ws.Range("K132").Value = 13357.9998
$ws.Range("L132").Value = 23792.7285
$ws.Range("M132").Value = -10827.9998
$ws.Range("N132").Value = -28852.7285

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6135.522
$ws.Range("I70").Value = 5237.091
$ws.Range("K70").Value = 5237.091
$ws.Range("M70").Value = -4967.091
$ws.Range("H73").Value = 6135.522
$ws.Range("I73").Value = 5237.091
$ws.Range("K73").Value = 5237.091
$ws.Range("M73").Value = -4301.091
$ws.Range("H122").Value = 3834.7046
$ws.Range("I122").Value = 3597.8684
$ws.Range("J122").Value = 5334.6665
$ws.Range("K122").Value = 10793.6052
$ws.Range("L122").Value = 16003.9995
$ws.Range("M122").Value = -8343.6052
$ws.Range("N122").Value = -20903.9995
$ws.Range("H126").Value = 3311.111
$ws.Range("I126").Value = 2433.3333
$ws.Range("J126").Value = 3750
$ws.Range("K126").Value = 7299.999899999999
$ws.Range("L126").Value = 11250
$ws.Range("M126").Value = -4829.999899999999
$ws.Range("N126").Value = -16190

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 15929.6
$ws.Range("I61").Value = 19078.834
$ws.Range("J61").Value = 3332.6667
$ws.Range("K61").Value = 19078.834
$ws.Range("L61").Value = 3332.6667
$ws.Range("M61").Value = -18876.834
$ws.Range("N61").Value = -3736.6667
$ws.Range("H82").Value = 2000
$ws.Range("J82").Value = 2000
$ws.Range("L82").Value = 2000
$ws.Range("N82").Value = -2722
$ws.Range("H85").Value = 2000
$ws.Range("J85").Value = 2000
$ws.Range("L85").Value = 2000
$ws.Range("N85").Value = -4496
$ws.Range("H98").Value = 26000
$ws.Range("J98").Value = 26000
$ws.Range("L98").Value = 26000
$ws.Range("N98").Value = -31990
$ws.Range("H113").Value = 15929.6
$ws.Range("I113").Value = 19078.834
$ws.Range("J113").Value = 3332.6667
$ws.Range("K113").Value = 19078.834
$ws.Range("L113").Value = 3332.6667
$ws.Range("M113").Value = -16908.834
$ws.Range("N113").Value = -7672.6667
$ws.Range("H122").Value = 1019547.1
$ws.Range("I122").Value = 1553482.2
$ws.Range("J122").Value = 5070.4
$ws.Range("K122").Value = 4660446.6
$ws.Range("L122").Value = 15211.2
$ws.Range("M122").Value = -4657996.6
$ws.Range("N122").Value = -20111.2
$ws.Range("H132").Value = 4907.7
$ws.Range("I132").Value = 4597.294
$ws.Range("K132").Value = 13791.882
$ws.Range("M132").Value = -11261.882

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 19704470
$ws.Range("I39").Value = 39393940
$ws.Range("J39").Value = 15000
$ws.Range("K39").Value = 39393940
$ws.Range("L39").Value = 15000
$ws.Range("M39").Value = -39393527
$ws.Range("N39").Value = -15826
$ws.Range("H42").Value = 23750
$ws.Range("H43").Value = 23750
$ws.Range("I43").Value = 50000
$ws.Range("J43").Value = 15000
$ws.Range("K43").Value = 50000
$ws.Range("L43").Value = 15000
$ws.Range("M43").Value = -49851
$ws.Range("N43").Value = -15298
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H113").Value = 569.26666
$ws.Range("I113").Value = 528.25
$ws.Range("J113").Value = 733.3333
$ws.Range("K113").Value = 1584.75
$ws.Range("L113").Value = 2199.9999
$ws.Range("M113").Value = 585.25
$ws.Range("N113").Value = -6539.9999
$ws.Range("H132").Value = 2367.6296
$ws.Range("I132").Value = 2029.8572
$ws.Range("K132").Value = 6089.571599999999
$ws.Range("M132").Value = -3559.571599999999
